$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44386
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 11000
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 11500
$ws.Range("P2").Value = 460

# Row 4
$ws.Range("D4").Value = 44370
$ws.Range("J4").Value = 520

# Row 5
$ws.Range("D5").Value = 44484
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 9000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 9500
$ws.Range("O5").Value = "Provincia del Elquí"
$ws.Range("P5").Value = 380

# Row 6
$ws.Range("D6").Value = 44376
$ws.Range("K6").Value = 12000
$ws.Range("L6").Value = 13000
$ws.Range("M6").Value = 12500
$ws.Range("P6").Value = 500

# Row 7
$ws.Range("D7").Value = 44446
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 11000
$ws.Range("M7").Value = 11500
$ws.Range("P7").Value = 460

# Row 8
$ws.Range("D8").Value = 44425
$ws.Range("K8").Value = 11500
$ws.Range("L8").Value = 12000
$ws.Range("M8").Value = 11750
$ws.Range("P8").Value = 470

# Row 9
$ws.Range("D9").Value = 44356
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 13000
$ws.Range("L9").Value = 14000
$ws.Range("M9").Value = 13500
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 540

# Row 10
$ws.Range("D10").Value = 44372
$ws.Range("J10").Value = 500

# Row 11
$ws.Range("D11").Value = 44316
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 16000
$ws.Range("L11").Value = 17000
$ws.Range("M11").Value = 16500
$ws.Range("P11").Value = 660

# Row 12
$ws.Range("D12").Value = 44473
$ws.Range("K12").Value = 8500
$ws.Range("L12").Value = 9000
$ws.Range("M12").Value = 8750
$ws.Range("P12").Value = 350

# Row 13
$ws.Range("D13").Value = 44466
$ws.Range("J13").Value = 400
$ws.Range("K13").Value = 9500
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = 9750
$ws.Range("P13").Value = 390

# Row 14
$ws.Range("D14").Value = 44384
$ws.Range("J14").Value = 560
$ws.Range("K14").Value = 11500
$ws.Range("L14").Value = 12000
$ws.Range("M14").Value = 11750
$ws.Range("P14").Value = 470
